$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '  +3.18%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.199.39'
$ws.Range("E3").Value = '  +0.97%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '260.30'
$ws.Range("E5").Value = '  +3.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '82.58'
$ws.Range("E6").Value = '  +12.66%  '
$ws.Range("E7").Value = '  +2.84%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +1.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.68'
$ws.Range("E10").Value = '  +8.87%  '
$ws.Range("E11").Value = '  +1.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.98'
$ws.Range("E12").Value = '  +3.27%  '
$ws.Range("E13").Value = '  +2.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.527.05'
$ws.Range("E14").Value = '  +0.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.30'
$ws.Range("E15").Value = '  +1.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.187.84'
$ws.Range("E16").Value = '  +0.11%  '
$ws.Range("E17").Value = '  +1.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.602.69'
$ws.Range("E18").Value = '  +3.23%  '
$ws.Range("E19").Value = '  +0.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.78'
$ws.Range("E20").Value = '  -1.18%  '
$ws.Range("E21").Value = '  +1.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.42'
$ws.Range("E22").Value = '  +15.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '230.81'
$ws.Range("E23").Value = '  +1.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.85'
$ws.Range("E24").Value = '  -4.96%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("B26").Value = 'InjectiveProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '42.53'
$ws.Range("E26").Value = '  +16.21%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.72'
$ws.Range("E27").Value = '  +2.34%  '
$ws.Range("B28").Value = 'WEMIXToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.37'
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.24'
$ws.Range("E29").Value = '  +3.73%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.20'
$ws.Range("E30").Value = '  +0.94%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '174.01'
$ws.Range("E31").Value = '  +2.36%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.43'
$ws.Range("E32").Value = '  +2.22%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0874'
$ws.Range("E33").Value = '  +6.81%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.32'
$ws.Range("E34").Value = '  +4.22%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.114'
$ws.Range("E35").Value = '  +6.66%  '
$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.122'
$ws.Range("E36").Value = '  +1.98%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.47'
$ws.Range("E37").Value = '  +6.62%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0352'
$ws.Range("E38").Value = '  +5.18%  '
$ws.Range("B39").Value = 'Celestia'
$ws.Range("C39").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.05'
$ws.Range("E39").Value = '  +11.00%  '
$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.86'
$ws.Range("E40").Value = '  +16.91%  '
$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.10'
$ws.Range("E41").Value = '  +2.47%  '
$ws.Range("B42").Value = 'MultiversX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '64.32'
$ws.Range("E42").Value = '  +8.77%  '
$ws.Range("B43").Value = 'THORChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.47'
$ws.Range("E43").Value = '  +6.52%  '
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.200'
$ws.Range("E44").Value = '  +2.51%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '100.43'
$ws.Range("E45").Value = '  -1.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0979'
$ws.Range("E46").Value = '  +0.98%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.24'
$ws.Range("E47").Value = '  +1.14%  '
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.17'
$ws.Range("E48").Value = '  +4.08%  '
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.11'
$ws.Range("E49").Value = '  +2.76%  '
$ws.Range("B50").Value = 'WOONetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.441'
$ws.Range("E50").Value = '  -3.66%  '
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.49'
$ws.Range("E51").Value = '  +23.20%  '
